$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 2-18 down to 3-19 (copy from bottom up to avoid overwrite)
for ($r = 18; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":E" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":E" + ($r + 1))
    $src.Copy($dst)
}

# Update dimension-relevant data: insert new row 2 values, and refresh C/E columns for all data rows
$ws.Range("A2").Value2 = 39400
$ws.Range("B2").Value2 = 2007
$ws.Range("C2").Value2 = -0.7015558851707571
$ws.Range("D2").Value2 = 2008
$ws.Range("E2").ClearContents()

$ws.Range("A3").Value2 = 39765
$ws.Range("B3").Value2 = 2008
$ws.Range("C3").Value2 = 0.3590181115727509
$ws.Range("D3").Value2 = 2009
$ws.Range("E3").ClearContents()

$ws.Range("A4").Value2 = 40130
$ws.Range("B4").Value2 = 2009
$ws.Range("C4").Value2 = -0.01587181126743165
$ws.Range("D4").Value2 = 2010
$ws.Range("E4").ClearContents()

$ws.Range("A5").Value2 = 40494
$ws.Range("B5").Value2 = 2010
$ws.Range("C5").Value2 = -0.02256889165885845
$ws.Range("D5").Value2 = 2011
$ws.Range("E5").Value2 = 0.1850158025574977

$ws.Range("A6").Value2 = 40862
$ws.Range("B6").Value2 = 2011
$ws.Range("C6").Value2 = 0.09611428386597787
$ws.Range("D6").Value2 = 2012
$ws.Range("E6").Value2 = -0.0461580488825053

$ws.Range("A7").Value2 = 41228
$ws.Range("B7").Value2 = 2012
$ws.Range("C7").Value2 = -0.1827723404408288
$ws.Range("D7").Value2 = 2013
$ws.Range("E7").Value2 = -0.1118170657869455

$ws.Range("A8").Value2 = 41592
$ws.Range("B8").Value2 = 2013
$ws.Range("C8").Value2 = -0.001350220946483294
$ws.Range("D8").Value2 = 2014
$ws.Range("E8").Value2 = -0.02003709364999384

$ws.Range("A9").Value2 = 41957
$ws.Range("B9").Value2 = 2014
$ws.Range("C9").Value2 = -0.075754880139145
$ws.Range("D9").Value2 = 2015
$ws.Range("E9").Value2 = -0.2383077634182995

$ws.Range("A10").Value2 = 42321
$ws.Range("B10").Value2 = 2015
$ws.Range("C10").Value2 = -0.5761528471665445
$ws.Range("D10").Value2 = 2016
$ws.Range("E10").Value2 = -0.2318455351884685

$ws.Range("A11").Value2 = 42689
$ws.Range("B11").Value2 = 2016
$ws.Range("C11").Value2 = -0.2011999787958185
$ws.Range("D11").Value2 = 2017
$ws.Range("E11").Value2 = -0.2555440101933648

$ws.Range("A12").Value2 = 43053
$ws.Range("B12").Value2 = 2017
$ws.Range("C12").Value2 = 0.1213692818849532
$ws.Range("D12").Value2 = 2018
$ws.Range("E12").Value2 = -0.06802695342326137

$ws.Range("A13").Value2 = 43418
$ws.Range("B13").Value2 = 2018
$ws.Range("C13").Value2 = 0.1493219406571766
$ws.Range("D13").Value2 = 2019
$ws.Range("E13").Value2 = 0.2006993856294326

$ws.Range("A14").Value2 = 43783
$ws.Range("B14").Value2 = 2019
$ws.Range("C14").Value2 = -0.4278219446121612
$ws.Range("D14").Value2 = 2020
$ws.Range("E14").Value2 = -0.5866687821558636

$ws.Range("A15").Value2 = 44159
$ws.Range("B15").Value2 = 2020
$ws.Range("C15").Value2 = -1.026566979837418
$ws.Range("D15").Value2 = 2021
$ws.Range("E15").Value2 = -2.497698913894009

$ws.Range("A16").Value2 = 44525
$ws.Range("B16").Value2 = 2021
$ws.Range("C16").Value2 = 0.3179894933462268
$ws.Range("D16").Value2 = 2022
$ws.Range("E16").Value2 = 0.1527842396192636

$ws.Range("A17").Value2 = 44890
$ws.Range("B17").Value2 = 2022
$ws.Range("C17").Value2 = 0.463604920919658
$ws.Range("D17").Value2 = 2023
$ws.Range("E17").Value2 = -0.2029074457040814

$ws.Range("A18").Value2 = 45254
$ws.Range("B18").Value2 = 2023
$ws.Range("C18").Value2 = 0.621639092134818
$ws.Range("D18").Value2 = 2024
$ws.Range("E18").Value2 = -0.1003146131857147

$ws.Range("A19").Value2 = 45618
$ws.Range("B19").Value2 = 2024
$ws.Range("C19").Value2 = -0.6768900623516982
$ws.Range("D19").Value2 = 2025
$ws.Range("E19").Value2 = 0.6778422458094902
